$wb = $excel.ActiveWorkbook

# --- Rename worksheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Role Schema"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Interaction Model"

# ====================================================================
# Sheet "Role Schema": insert two new role-schema blocks before the
# first existing block (which currently occupies rows 4-9).
# Each block is 8 rows tall (6 data rows + 2 separator rows), so we
# insert 16 rows at row 4 to make room for two new blocks.
# ====================================================================
$ws1.Range("F4:F19").EntireRow.Insert() | Out-Null

# Restore the blank separator rows' height (rows 11 and 19), matching
# the other blank separator rows in the sheet (e.g. row 3).
$ws1.Rows.Item(11).RowHeight = 15
$ws1.Rows.Item(19).RowHeight = 15

# Copy formatting (styles/borders/merges) from the block that is now at
# rows 20-25 (the original first block) into the two new blank blocks.
$ws1.Range("F20:H25").Copy() | Out-Null
$ws1.Range("F4").PasteSpecial(-4122) | Out-Null
$ws1.Range("F20:H25").Copy() | Out-Null
$ws1.Range("F12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Match the row heights used by the rest of the data rows (the new
# blocks should look the same as the other role-schema blocks).
for ($r = 4; $r -le 9; $r++) { $ws1.Rows.Item($r).RowHeight = 15 }
for ($r = 12; $r -le 17; $r++) { $ws1.Rows.Item($r).RowHeight = 15 }

# Re-create the merges for the two new blocks (PasteSpecial of formats
# alone does not recreate merged ranges).
$ws1.Range("F4:G4").Merge() | Out-Null
$ws1.Range("F5:G5").Merge() | Out-Null
$ws1.Range("F6:G6").Merge() | Out-Null
$ws1.Range("F7:G7").Merge() | Out-Null
$ws1.Range("F8:F9").Merge() | Out-Null

$ws1.Range("F12:G12").Merge() | Out-Null
$ws1.Range("F13:G13").Merge() | Out-Null
$ws1.Range("F14:G14").Merge() | Out-Null
$ws1.Range("F15:G15").Merge() | Out-Null
$ws1.Range("F16:F17").Merge() | Out-Null

# --- Block 1 (rows 4-9): "searchProfile" role ---
$ws1.Range("F4").Value = "Role Schema"
$ws1.Range("H4").Value = "searchProfile"

$ws1.Range("F5").Value = "Description"
$ws1.Range("H5").Value = "Search in Database for user profile"

$ws1.Range("F6").Value = "Protocols and Activities"
$ws1.Range("H6").Value = "SearchDatabase"

$ws1.Range("F7").Value = "Permissions"
$ws1.Range("H7").Value = "read Userprofile, write Database"

$ws1.Range("F8").Value = "Responsibilities"
$ws1.Range("G8").Value = "Liveness"
$ws1.Range("H8").Value = "newProfile = SearchDatabase.service"

$ws1.Range("G9").Value = "Safety"
$ws1.Range("H9").Value = "Successfully creating a profile"

# --- Block 2 (rows 12-17): "createAppointment" role ---
$ws1.Range("F12").Value = "Role Schema"
$ws1.Range("H12").Value = "createAppointment"

$ws1.Range("F13").Value = "Description"
$ws1.Range("H13").Value = "Creating new and successful appointment"

$ws1.Range("F14").Value = "Protocols and Activities"
$ws1.Range("H14").Value = "createAppointment"

$ws1.Range("F15").Value = "Permissions"
$ws1.Range("H15").Value = "read Userprofile, write Database"

$ws1.Range("F16").Value = "Responsibilities"
$ws1.Range("G16").Value = "Liveness"
$ws1.Range("H16").Value = "createAppointment = createApppointment.service"

$ws1.Range("G17").Value = "Safety"
$ws1.Range("H17").Value = "Successfully generating an Appointment"

# ====================================================================
# Sheet "Interaction Model": fill in the previously blank E and F
# columns with the new "searchDatabase" / "createAppointment" agents.
# ====================================================================
$ws2.Range("E5").Value = "searchDatabase"
$ws2.Range("F5").Value = "createAppointment"

$ws2.Range("E6").Value = "To search in databases for any user profile already exists"
$ws2.Range("F6").Value = "To create a new and repeated appointments with the doctor for the users  "

$ws2.Range("E7").Value = "Clinic system"
$ws2.Range("F7").Value = "Appointment"

$ws2.Range("E8").Value = "searchDatabse"
$ws2.Range("F8").Value = "createAppointment"

$ws2.Range("E9").Value = "If any registered user visits portal, it will search the database for its previous logs and details and share it with the physician. "
$ws2.Range("F9").Value = "When a user wish to consult a physician, this service will generate appointments for the user with the physician. It will store the appointments in the database so that reminders can be sent."
